# Generate Report for Handoff
#
# A new handoff/generate run replaced the old file-id
# "e9d340a2-ca04-4764-8c07-824497279858" with a new one
# "3f812db5-a0b6-4342-b9ef-b92334164cfb" (and a new content hash for the
# generated .xlf files), and the handoff timestamps moved forward a bit.
# Update every cell (and hyperlink display text) that mirrors those values.

$wb = $excel.ActiveWorkbook

$newId = "3f812db5-a0b6-4342-b9ef-b92334164cfb"
$newHash = "5e445a842f46d5866b277a342dd65f5f27f5a028"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
$wsOverview.Range("G2").Value = "2016-08-30 03:01:14"

foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = "e2e\$newId.md"
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newId.md"
$wsZhCn.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-30 03:01:08"

foreach ($h in $wsZhCn.Hyperlinks) {
    $h.TextToDisplay = "$newId.md"
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newId.md"
$wsDeDe.Range("G2").Value = "$newId.$newHash.de-de.xlf"
# H2 ("Latest Handoff Datetime") shared its value with Overview!G2 in the
# original shared-string table, so it picks up the same new timestamp.
$wsDeDe.Range("H2").Value = "2016-08-30 03:01:14"

foreach ($h in $wsDeDe.Hyperlinks) {
    $h.TextToDisplay = "$newId.md"
}
